$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate() | Out-Null

# Insert a new blank column before column N (14), shifting N:P -> O:Q
$ws.Columns("N:N").Insert() | Out-Null
$ws.Columns("N:N").ColumnWidth = 9.166666666666666

# Update the selection on the Repayment Schedule sheet to match the new layout
$ws.Range("S8").Select() | Out-Null

# Restore the previous selection on the Transactions sheet
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("D14").Select() | Out-Null

# Re-activate Repayment Schedule so it is the tab shown when the workbook is opened
$ws.Activate() | Out-Null
